$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update B6 value from 2 to 4 (this drives the recalculation of D5 and F5)
$ws.Range("B6").Value = 4

# Move the active selection to B7 (matches the view state change in the diff)
$ws.Range("B7").Select()
